# Generate Report for Handoff
# Adds two new file rows (3d67bfa0-...  and 637254d8-...) to the
# "Overview", "zh-cn" and "de-de" worksheets of the handback-status report.

$wb = $excel.ActiveWorkbook

$newFile1 = "3d67bfa0-762f-44da-91ef-c03df45add61"
$newFile2 = "637254d8-c43c-442b-93fd-7ec1a1a12849"

$hash1 = "9c0cbf2c663cfc2e3ed19b428b84a84da312fa62"
$hash2 = "e02ca436e165f99c205235f56d5425e38295b319"

$status = "Ready for handoff"

$handoffDateOverview1 = "2016-03-24 09:59:13"   # de-de / overview handoff datetime for file1
$handoffDateOverview2 = "2016-03-24 09:59:13"   # de-de / overview handoff datetime for file2

$handoffDatetimeZh = "2016-03-24 09:59:04"      # zh-cn handoff datetime (both files)

$targetFileDate = "0001-01-01 00:00:00"
$handoffReason = "Include"
$fileExt = ".md"

$mdRepoBase = "https://github.com/OpenLocalizationTest/oltest/blob/2f6fd0799dbba112e420c2bc7ac4add691ee6923/e2e/"
$zhRepoBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9863f660a99afd9efd5e889bab3f4cc571bd5cf2/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/"
$deRepoBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/66189f2d89b15a6b9a78c193a4075a1552a9308f/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/"

# Pre-compute every composite string up front - the interpreter used here does
# not reliably evaluate parenthesised expressions passed inline as function
# arguments, so every value is built into its own variable first.

$md1 = $newFile1 + ".md"
$md2 = $newFile2 + ".md"

$mdUrl1 = $mdRepoBase + $md1
$mdUrl2 = $mdRepoBase + $md2

$zhXlf1 = $newFile1 + "." + $hash1 + ".zh-cn.xlf"
$zhXlf2 = $newFile2 + "." + $hash2 + ".zh-cn.xlf"
$deXlf1 = $newFile1 + "." + $hash1 + ".de-de.xlf"
$deXlf2 = $newFile2 + "." + $hash2 + ".de-de.xlf"

$zhXlfUrl1 = $zhRepoBase + $zhXlf1
$zhXlfUrl2 = $zhRepoBase + $zhXlf2
$deXlfUrl1 = $deRepoBase + $deXlf1
$deXlfUrl2 = $deRepoBase + $deXlf2

function Set-HyperlinkCell {
    param($ws, $cellAddr, $text, $url, $display)
    $ws.Range($cellAddr).Value = $text
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $url, [Type]::Missing, [Type]::Missing, $display) | Out-Null
    # Match the pre-existing "HyperLink" look (underline + cornflowerblue) used
    # elsewhere in the workbook for link cells.
    $ws.Range($cellAddr).Font.Underline = $true
    $ws.Range($cellAddr).Font.Color = 15570276   # BGR for RGB 6495ED (cornflowerblue)
}

function Set-DateCell {
    param($ws, $cellAddr, $text)
    $ws.Range($cellAddr).Value = $text
    $ws.Range($cellAddr).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 4 - 3d67bfa0-...
Set-HyperlinkCell $wsOverview "A4" $md1 $mdUrl1 $md1
$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status
Set-DateCell $wsOverview "D4" $handoffDateOverview1

# Row 5 - 637254d8-...
Set-HyperlinkCell $wsOverview "A5" $md2 $mdUrl2 $md2
$wsOverview.Range("B5").Value = $status
$wsOverview.Range("C5").Value = $status
Set-DateCell $wsOverview "D5" $handoffDateOverview2

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 4 - 3d67bfa0-...
Set-HyperlinkCell $wsZh "A4" $md1 $mdUrl1 $md1
$wsZh.Range("B4").Value = $fileExt
$wsZh.Range("C4").Value = $status
Set-HyperlinkCell $wsZh "D4" $zhXlf1 $zhXlfUrl1 $zhXlf1
Set-DateCell $wsZh "E4" $handoffDatetimeZh
Set-DateCell $wsZh "H4" $targetFileDate
$wsZh.Range("J4").Value = $handoffReason

# Row 5 - 637254d8-...
Set-HyperlinkCell $wsZh "A5" $md2 $mdUrl2 $md2
$wsZh.Range("B5").Value = $fileExt
$wsZh.Range("C5").Value = $status
Set-HyperlinkCell $wsZh "D5" $zhXlf2 $zhXlfUrl2 $zhXlf2
Set-DateCell $wsZh "E5" $handoffDatetimeZh
Set-DateCell $wsZh "H5" $targetFileDate
$wsZh.Range("J5").Value = $handoffReason

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 4 - 3d67bfa0-...
Set-HyperlinkCell $wsDe "A4" $md1 $mdUrl1 $md1
$wsDe.Range("B4").Value = $fileExt
$wsDe.Range("C4").Value = $status
Set-HyperlinkCell $wsDe "D4" $deXlf1 $deXlfUrl1 $deXlf1
Set-DateCell $wsDe "E4" $handoffDateOverview1
Set-DateCell $wsDe "H4" $targetFileDate
$wsDe.Range("J4").Value = $handoffReason

# Row 5 - 637254d8-...
Set-HyperlinkCell $wsDe "A5" $md2 $mdUrl2 $md2
$wsDe.Range("B5").Value = $fileExt
$wsDe.Range("C5").Value = $status
Set-HyperlinkCell $wsDe "D5" $deXlf2 $deXlfUrl2 $deXlf2
Set-DateCell $wsDe "E5" $handoffDateOverview2
Set-DateCell $wsDe "H5" $targetFileDate
$wsDe.Range("J5").Value = $handoffReason
